# Generate Report for Archive
# - Status has moved on from "Ready for handoff" to "In Translation" for the
#   two e2e test files, on every sheet that surfaces that status.
# - The Status column(s) are re-fit to the new, shorter text.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# Narrow the now-shorter Status columns to match the new content width.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
